$d = $word.ActiveDocument

# Find the paragraph "Install relevant software." and grab its range.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Install relevant software.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph 'Install relevant software.'"
}

# Collapse to the end of that paragraph (after the text, before the
# paragraph mark) and insert a brand new paragraph after it, inheriting
# the same list formatting (ListParagraph / numId 3).
$insertionPoint = $target.Range.Duplicate
$insertionPoint.Collapse(0)  # wdCollapseEnd
$newPara = $insertionPoint.InsertParagraphAfter()

# Move into the freshly created paragraph and type the new bullet text.
$newRange = $d.Range($target.Range.End, $target.Range.End)
$newRange.InsertAfter("Consult relevan")
